$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the first two detail rows (LUDMILA / RENATO) ---
# Row 2 = 000678704 / LUDMILA / 184000
# Row 3 (after the first delete) = 000330949 / RENATO / 50000
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# --- ALESSANDRA's balance is corrected from 30001 to 30011 ---
# (now sitting at row 2 after the two deletions above)
$ws.Cells.Item(2, 3).Value = 30011

# --- A new account (HIROKO) is inserted right after ALESSANDRA ---
$ws.Rows.Item(3).Insert()
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "004996634"
$ws.Cells.Item(3, 2).Value = "HIROKO"
$ws.Cells.Item(3, 3).Value = 27298.15

# --- POLYANNA / DIOGO / LOHRAN (now rows 5,5,5) are removed ... ---
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()

# --- ... and replaced by a single new account (JOAQUIM) ---
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "004321016"
$ws.Cells.Item(5, 2).Value = "JOAQUIM"
$ws.Cells.Item(5, 3).Value = 4052.1

# --- The old HIROKO row (-287.81) near the end of the list is removed ---
# (it sat at row 244 originally; net of the edits above (-2-3+1+1) it is
#  now 3 rows earlier, i.e. row 241)
$ws.Rows.Item(241).Delete()
